# Commit: "added charging mechanism based on parsed information"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OCS Input")

# Replace the placeholder "Charging Mechanism" text (column L) for each
# line item with the newly parsed charging-mechanism description.
$ws.Range("L11").Value = "1 unit/day from start phase 50 to end phase 50 for maximum 7 occurrences"
$ws.Range("L12").Value = "1 unit/day from start phase 15 to end phase 15 for maximum 7 occurrences"
$ws.Range("L13").Value = "1 unit/day from start phase 5 to end phase 5 for maximum 7 occurrences"
$ws.Range("L14").Value = "1 unit/day from start phase 20 to end phase 20 for maximum 7 occurrences"

# Correct the supplier / well code typo: GK-W527B -> GK-P527B
$ws.Range("B5").Value = "GK-P527B"

# Leave the selection where the user was last working
$ws.Range("B5").Select()
